$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.286.81"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "1.821.83"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.65"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4681"
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3774"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07414"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8724"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.60"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.822.21"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.677"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.410"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.45"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07109"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008788"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "27.293.86"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.310"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "2.046.53"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.41"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.251"
$ws.Range("E27").Value = "  +4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.307"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.28"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08932"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7803"
$ws.Range("E32").Value = "  +5.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.181"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.527"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.945"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.099"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01968"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05249"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.239"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5324"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.348"
$ws.Range("E43").Value = "  +21.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1693"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.616"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5076"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.34"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.674"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06347"
$ws.Range("E51").Value = "  +0.98%  "
